# Deploying to gh-pages — add the 2021 column (Y) to the undernourishment
# table and update the sheet view's selection, matching the upstream
# 2.1.1.xlsx commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for the 2021 column, keyed by row number.
$values = @{
    4  = 2021
    5  = 46.69
    6  = 52.52
    7  = 43.22
    8  = 51.31
    9  = 41.31
    10 = 52.43
    11 = 49.27
    12 = 31.68
    13 = 35.59
    14 = 55.28
    15 = 61.02
    16 = 48.72
}

# Column Y mirrors column X's formatting for every populated row (4-16),
# so copy each X cell's format/style into the new Y cell before writing
# the 2021 figure into it.
foreach ($row in 4..16) {
    $srcCell = $ws.Cells.Item($row, 24)   # column X
    $dstCell = $ws.Cells.Item($row, 25)   # column Y
    $srcCell.Copy($dstCell)
    $dstCell.Value = $values[$row]
}

# Scroll the view so column B is the left-most visible column (best effort;
# mirrors the authored sheetView's topLeftCell="B1"), then land the
# selection on AA15 as in the source workbook.
$excel.ActiveWindow.ScrollColumn = 2
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("AA15").Select()
